$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new row at position 359 - this pushes the existing rows
# 359..453 down to 360..454 (carrying all of their original values with
# them), growing the used range to A1:R454.
$ws.Rows.Item(359).Insert()

# Populate the newly inserted (blank) row 359 with the new weekly record.
$ws.Range("A359").Value = 5
$ws.Range("B359").Value = "Macroferia Regional de Talca"
$ws.Range("C359").Value = "Maule"
$ws.Range("D359").Value = 44932
$ws.Range("E359").Value = 7
$ws.Range("F359").Value = 100114014
$ws.Range("G359").Value = "Betarraga"
$ws.Range("H359").Value = "Sin especificar"
$ws.Range("I359").Value = "Primera"
$ws.Range("J359").Value = 5000
$ws.Range("K359").Value = 600
$ws.Range("L359").Value = 600
$ws.Range("M359").Value = 600
$ws.Range("N359").Value = "$/paquete 5 unidades"
$ws.Range("O359").Value = "Región del Maule"
$ws.Range("P359").Value = 120
$ws.Range("Q359").Value = 5
$ws.Range("R359").Value = "Hortaliza"
